# Add new Lincoln data
# Rename all "PalmerstonNorth1994..." simulation names (column A) to "Lincoln1994..."
# and update the saved worksheet view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "PalmerstonNorth1994*") {
        $cell.Value = $val -replace "^PalmerstonNorth1994", "Lincoln1994"
    }
}

# Update the view state to match the saved window position/selection.
$window = $excel.ActiveWindow
$window.ScrollRow = 221
$window.ScrollColumn = 1
$ws.Range("A133").Select()
